$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row strings (row 1)
$ws.Range("C1").Value = "A_C"
$ws.Range("D1").Value = "FFR_A"
$ws.Range("E1").Value = "C_A"

# Update data values in row 2 and row 3 for columns B:E
$ws.Range("B2").Value = -78.29729548487525
$ws.Range("C2").Value = 10.19294504512097
$ws.Range("D2").Value = -1.821719970535869
$ws.Range("E2").Value = 0.08755023134145418

$ws.Range("B3").Value = 0.000001288281721301132
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.000000000001849409514420586
$ws.Range("E3").Value = 0

# Remove now-unused columns F and G entirely
$ws.Range("F1:G3").Delete()
